$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2953.6667
$ws.Range("I28").Value = 3833
$ws.Range("J28").Value = 2074.3333
$ws.Range("K28").Value = 3833
$ws.Range("L28").Value = 2074.3333
$ws.Range("M28").Value = -3348
$ws.Range("N28").Value = -3044.3333
$ws.Range("H32").Value = 5619.3335
$ws.Range("I32").Value = 3477.7778
$ws.Range("J32").Value = 8831.666999999999
$ws.Range("K32").Value = 3477.7778
$ws.Range("L32").Value = 8831.666999999999
$ws.Range("M32").Value = -3151.7778
$ws.Range("N32").Value = -9483.666999999999
$ws.Range("H51").Value = 11827.357
$ws.Range("I51").Value = 34249.5
$ws.Range("J51").Value = 8090.3335
$ws.Range("K51").Value = 34249.5
$ws.Range("L51").Value = 8090.3335
$ws.Range("M51").Value = -33765.5
$ws.Range("N51").Value = -9058.333500000001
$ws.Range("H87").Value = 51998.4
$ws.Range("J87").Value = 51998.4
$ws.Range("L87").Value = 51998.4
$ws.Range("N87").Value = -54494.4
$ws.Range("H90").Value = 51998.4
$ws.Range("J90").Value = 51998.4
$ws.Range("L90").Value = 155995.2
$ws.Range("N90").Value = -168475.2
$ws.Range("H92").Value = 1527.625
$ws.Range("I92").Value = 946.2
$ws.Range("K92").Value = 946.2
$ws.Range("M92").Value = 301.8
$ws.Range("H135").Value = 476726.28
$ws.Range("I135").Value = 476726.28
$ws.Range("K135").Value = 4290536.52
$ws.Range("M135").Value = -4288001.52
$ws.Range("H137").Value = 2351.963
$ws.Range("I137").Value = 2607.8823
$ws.Range("K137").Value = 7823.646900000001
$ws.Range("M137").Value = -5273.646900000001
$ws.Range("H138").Value = 2470.2856
$ws.Range("J138").Value = 2525
$ws.Range("L138").Value = 7575
$ws.Range("N138").Value = -17855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2337.3333
$ws.Range("I2").Value = 804.4
$ws.Range("K2").Value = 804.4
$ws.Range("M2").Value = -691.4
$ws.Range("H5").Value = 267.1111
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
$ws.Range("H32").Value = 5072.4844
$ws.Range("I32").Value = 5020.783
$ws.Range("K32").Value = 5020.783
$ws.Range("M32").Value = -4733.783
$ws.Range("H116").Value = 2337.3333
$ws.Range("I116").Value = 804.4
$ws.Range("K116").Value = 804.4
$ws.Range("M116").Value = 1489.6
$ws.Range("H132").Value = 3040.6047
$ws.Range("J132").Value = 5972.923
$ws.Range("L132").Value = 17918.769
$ws.Range("N132").Value = -22978.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2337.3333
$ws.Range("I3").Value = 804.4
$ws.Range("K3").Value = 804.4
$ws.Range("M3").Value = -690.4
$ws.Range("H4").Value = 267.1111
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = ""
$ws.Range("N51").Value = 0
$ws.Range("H88").Value = 40343
$ws.Range("J88").Value = 40343
$ws.Range("L88").Value = 40343
$ws.Range("N88").Value = -41155
$ws.Range("H91").Value = 40343
$ws.Range("J91").Value = 40343
$ws.Range("L91").Value = 40343
$ws.Range("N91").Value = -43151
$ws.Range("H122").Value = 34319.6
$ws.Range("J122").Value = 34319.6
$ws.Range("L122").Value = 34319.6
$ws.Range("N122").Value = -44119.6
$ws.Range("H134").Value = 5281.6772
$ws.Range("I134").Value = 2441.25
$ws.Range("K134").Value = 7323.75
$ws.Range("M134").Value = -4788.75
$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8396.259
$ws.Range("I58").Value = 2543.125
$ws.Range("K58").Value = 2543.125
$ws.Range("M58").Value = -2340.125
$ws.Range("H59").Value = 66999.8
$ws.Range("J59").Value = 88333
$ws.Range("L59").Value = 88333
$ws.Range("N59").Value = -90623
$ws.Range("H60").Value = 16856.428
$ws.Range("J60").Value = 37497.5
$ws.Range("L60").Value = 37497.5
$ws.Range("N60").Value = -38519.5
$ws.Range("H132").Value = 10261804
$ws.Range("I132").Value = 2387.4348
$ws.Range("K132").Value = 7162.3044
$ws.Range("M132").Value = -4632.3044
$ws.Range("H134").Value = 7100.7144
$ws.Range("I134").Value = 1553.0834
$ws.Range("J134").Value = 11261.4375
$ws.Range("K134").Value = 4659.2502
$ws.Range("L134").Value = 33784.3125
$ws.Range("M134").Value = -2124.2502
$ws.Range("N134").Value = -38854.3125
$ws.Range("H136").Value = 8396.259
$ws.Range("I136").Value = 2543.125
$ws.Range("K136").Value = 7629.375
$ws.Range("M136").Value = -5079.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3552.5
$ws.Range("I58").Value = 3552.5
$ws.Range("K58").Value = 10657.5
$ws.Range("M58").Value = -10529.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4158.5
$ws.Range("I80").Value = 3814
$ws.Range("J80").Value = 4365.2
$ws.Range("K80").Value = 3814
$ws.Range("L80").Value = 4365.2
$ws.Range("M80").Value = -2816
$ws.Range("N80").Value = -6361.2
$ws.Range("H83").Value = 4158.5
$ws.Range("I83").Value = 3814
$ws.Range("J83").Value = 4365.2
$ws.Range("K83").Value = 19070
$ws.Range("L83").Value = 21826
$ws.Range("M83").Value = -14078
$ws.Range("N83").Value = -31810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1126.1765
$ws.Range("I22").Value = 664.5925999999999
$ws.Range("J22").Value = 2906.5715
$ws.Range("K22").Value = 664.5925999999999
$ws.Range("L22").Value = 2906.5715
$ws.Range("M22").Value = -369.5925999999999
$ws.Range("N22").Value = -3496.5715
$ws.Range("H27").Value = 1126.1765
$ws.Range("I27").Value = 664.5925999999999
$ws.Range("J27").Value = 2906.5715
$ws.Range("K27").Value = 664.5925999999999
$ws.Range("L27").Value = 2906.5715
$ws.Range("M27").Value = -557.5925999999999
$ws.Range("N27").Value = -3120.5715
$ws.Range("H132").Value = 4933
$ws.Range("I132").Value = 2961.3928
$ws.Range("J132").Value = 7999.9443
$ws.Range("K132").Value = 8884.178400000001
$ws.Range("L132").Value = 23999.8329
$ws.Range("M132").Value = -6354.178400000001
$ws.Range("N132").Value = -29059.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10005333
$ws.Range("J81").Value = 40013020
$ws.Range("L81").Value = 80026040
$ws.Range("N81").Value = -80028162
$ws.Range("H84").Value = 10005333
$ws.Range("J84").Value = 40013020
$ws.Range("L84").Value = 400130200
$ws.Range("N84").Value = -400140808
$ws.Range("H126").Value = 2399.7
$ws.Range("I126").Value = 999.75
$ws.Range("K126").Value = 2999.25
$ws.Range("M126").Value = -529.25
$ws.Range("H132").Value = 13833.767
$ws.Range("I132").Value = 7335.5293
$ws.Range("K132").Value = 22006.5879
$ws.Range("M132").Value = -19476.5879
